$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old row 39 so the existing weekly records
# (previously rows 39-49) shift down to rows 41-51, preserving their data
# and formatting (the date-formatted style on column D follows the shift).
$ws.Rows("39:40").Insert()

# New week's data: row 39 - "Primera" quality
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 44504
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 300000000
$ws.Range("G39").Value = "Espárragos"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 800
$ws.Range("K39").Value = 1200
$ws.Range("L39").Value = 1300
$ws.Range("M39").Value = 1244
$ws.Range("N39").Value = "$/kilo"
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 1244
$ws.Range("Q39").Value = 1
$ws.Range("R39").Value = "Hortaliza"

# New week's data: row 40 - "Segunda" quality
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "Vega Modelo de Temuco"
$ws.Range("C40").Value = "La Araucanía"
$ws.Range("D40").Value = 44504
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = 300000000
$ws.Range("G40").Value = "Espárragos"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Segunda"
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = 1000
$ws.Range("N40").Value = "$/kilo"
$ws.Range("O40").Value = "Región del Maule"
$ws.Range("P40").Value = 1000
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = "Hortaliza"
